$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 117.85714
$ws.Range("I5").Value = 106
$ws.Range("J5").Value = 147.5
$ws.Range("K5").Value = 106
$ws.Range("L5").Value = 147.5
$ws.Range("M5").Value = 9
$ws.Range("N5").Value = -377.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 16100
$ws.Range("I18").Value = 20150
$ws.Range("J18").Value = 8000
$ws.Range("K18").Value = 20150
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = -19866
$ws.Range("N18").Value = -8568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5104.8125
$ws.Range("I32").Value = 3154.25
$ws.Range("J32").Value = 5755
$ws.Range("K32").Value = 3154.25
$ws.Range("L32").Value = 5755
$ws.Range("M32").Value = -2828.25
$ws.Range("N32").Value = -6407

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 864.125
$ws.Range("I33").Value = 187.08333
$ws.Range("K33").Value = 187.08333
$ws.Range("M33").Value = 41.91667000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4999.5
$ws.Range("J40").Value = 4999.4287
$ws.Range("L40").Value = 4999.4287
$ws.Range("N40").Value = -5349.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 45099.25
$ws.Range("I46").Value = 45099.25
$ws.Range("K46").Value = 135297.75
$ws.Range("M46").Value = -135178.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 122.052635
$ws.Range("I55").Value = 47.444443
$ws.Range("J55").Value = 189.2
$ws.Range("K55").Value = 47.444443
$ws.Range("L55").Value = 189.2
$ws.Range("M55").Value = 166.555557
$ws.Range("N55").Value = -617.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 45099.25
$ws.Range("I60").Value = 45099.25
$ws.Range("K60").Value = 135297.75
$ws.Range("M60").Value = -134813.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 200006400
$ws.Range("I69").Value = 6000
$ws.Range("K69").Value = 18000
$ws.Range("M69").Value = -17126

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5886865
$ws.Range("I70").Value = 12502858
$ws.Range("J70").Value = 5982.722
$ws.Range("K70").Value = 37508574
$ws.Range("L70").Value = 17948.166
$ws.Range("M70").Value = -37508304
$ws.Range("N70").Value = -18488.166

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 200006400
$ws.Range("I72").Value = 6000
$ws.Range("K72").Value = 54000
$ws.Range("M72").Value = -49632

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5886865
$ws.Range("I73").Value = 12502858
$ws.Range("J73").Value = 5982.722
$ws.Range("K73").Value = 37508574
$ws.Range("L73").Value = 17948.166
$ws.Range("M73").Value = -37507638
$ws.Range("N73").Value = -19820.166

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 681.4286
$ws.Range("I80").Value = 773.5
$ws.Range("J80").Value = 644.6
$ws.Range("K80").Value = 2320.5
$ws.Range("L80").Value = 1933.8
$ws.Range("M80").Value = -1322.5
$ws.Range("N80").Value = -3929.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 681.4286
$ws.Range("I83").Value = 773.5
$ws.Range("J83").Value = 644.6
$ws.Range("K83").Value = 6961.5
$ws.Range("L83").Value = 5801.400000000001
$ws.Range("M83").Value = -1969.5
$ws.Range("N83").Value = -15785.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1831.3077
$ws.Range("I86").Value = 1327.8334
$ws.Range("J86").Value = 2262.8572
$ws.Range("K86").Value = 1327.8334
$ws.Range("L86").Value = 2262.8572
$ws.Range("M86").Value = -204.8334
$ws.Range("N86").Value = -4508.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 6359.8667
$ws.Range("J88").Value = 7192.231
$ws.Range("L88").Value = 7192.231
$ws.Range("N88").Value = -8004.231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1831.3077
$ws.Range("I89").Value = 1327.8334
$ws.Range("J89").Value = 2262.8572
$ws.Range("K89").Value = 6639.166999999999
$ws.Range("L89").Value = 11314.286
$ws.Range("M89").Value = -1023.166999999999
$ws.Range("N89").Value = -22546.286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 6359.8667
$ws.Range("J91").Value = 7192.231
$ws.Range("L91").Value = 7192.231
$ws.Range("N91").Value = -10000.231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1425.6923
$ws.Range("I98").Value = 1434.72
$ws.Range("K98").Value = 1434.72
$ws.Range("M98").Value = 63.27999999999997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3013.2903
$ws.Range("I100").Value = 1744.9333
$ws.Range("J100").Value = 4202.375
$ws.Range("K100").Value = 1744.9333
$ws.Range("L100").Value = 4202.375
$ws.Range("M100").Value = -1203.9333
$ws.Range("N100").Value = -5284.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1193.421
$ws.Range("I103").Value = 1154.3572
$ws.Range("J103").Value = 1216.2084
$ws.Range("K103").Value = 3463.0716
$ws.Range("L103").Value = 3648.6252
$ws.Range("M103").Value = -2877.0716
$ws.Range("N103").Value = -4820.6252

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 8101.7617
$ws.Range("I107").Value = 8101.7617
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 8101.7617
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -6181.7617
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1425.6923
$ws.Range("I122").Value = 1434.72
$ws.Range("K122").Value = 4304.16
$ws.Range("M122").Value = -1854.16

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 16333.5
$ws.Range("J125").Value = 30199.715
$ws.Range("L125").Value = 271797.435
$ws.Range("N125").Value = -276717.435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1239.4546
$ws.Range("I127").Value = 1239.4546
$ws.Range("K127").Value = 3718.3638
$ws.Range("M127").Value = 1241.6362

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1695.0834
$ws.Range("I129").Value = 693.75
$ws.Range("J129").Value = 2195.75
$ws.Range("K129").Value = 2081.25
$ws.Range("L129").Value = 6587.25
$ws.Range("M129").Value = 2918.75
$ws.Range("N129").Value = -16587.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 99940
$ws.Range("J133").Value = 99940
$ws.Range("L133").Value = 99940
$ws.Range("N133").Value = -110060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9302.317999999999
$ws.Range("I137").Value = 13297.692
$ws.Range("J137").Value = 3531.2222
$ws.Range("K137").Value = 39893.076
$ws.Range("L137").Value = 10593.6666
$ws.Range("M137").Value = -37343.076
$ws.Range("N137").Value = -15693.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3374.587
$ws.Range("J138").Value = 3928
$ws.Range("L138").Value = 11784
$ws.Range("N138").Value = -22064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 30020.5
$ws.Range("J47").Value = 30020.5
$ws.Range("L47").Value = 30020.5
$ws.Range("N47").Value = -31470.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6315.154
$ws.Range("I61").Value = 6633.0835
$ws.Range("K61").Value = 6633.0835
$ws.Range("M61").Value = -6421.0835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2605
$ws.Range("I74").Value = 3062.2727
$ws.Range("K74").Value = 3062.2727
$ws.Range("M74").Value = -2188.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2605
$ws.Range("I77").Value = 3062.2727
$ws.Range("K77").Value = 15311.3635
$ws.Range("M77").Value = -10943.3635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3082.9333
$ws.Range("I88").Value = 3200
$ws.Range("K88").Value = 3200
$ws.Range("M88").Value = -2794

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3082.9333
$ws.Range("I91").Value = 3200
$ws.Range("K91").Value = 3200
$ws.Range("M91").Value = -1796

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1305.8462
$ws.Range("I122").Value = 1361.5454
$ws.Range("K122").Value = 4084.6362
$ws.Range("M122").Value = -1634.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6315.154
$ws.Range("I136").Value = 6633.0835
$ws.Range("K136").Value = 19899.2505
$ws.Range("M136").Value = -17349.2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1098.3529
$ws.Range("I20").Value = 1015.7273
$ws.Range("J20").Value = 1249.8334
$ws.Range("K20").Value = 1015.7273
$ws.Range("L20").Value = 1249.8334
$ws.Range("M20").Value = -768.7273
$ws.Range("N20").Value = -1743.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 890.9
$ws.Range("I22").Value = 879.8889
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 879.8889
$ws.Range("L22").Value = 990
$ws.Range("M22").Value = -706.8889
$ws.Range("N22").Value = -1336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 89666.664
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 299998.5
$ws.Range("J70").Value = 299998.5
$ws.Range("L70").Value = 299998.5
$ws.Range("N70").Value = -300584.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H73").Value = 299998.5
$ws.Range("J73").Value = 299998.5
$ws.Range("L73").Value = 299998.5
$ws.Range("N73").Value = -302026.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3030.6
$ws.Range("I86").Value = 3037.5
$ws.Range("K86").Value = 3037.5
$ws.Range("M86").Value = -1914.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3030.6
$ws.Range("I89").Value = 3037.5
$ws.Range("K89").Value = 15187.5
$ws.Range("M89").Value = -9571.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 94635.05499999999
$ws.Range("I94").Value = 151772.61
$ws.Range("K94").Value = 151772.61
$ws.Range("M94").Value = -151321.61

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6002.2
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6742.9
$ws.Range("I107").Value = 1857.25
$ws.Range("K107").Value = 1857.25
$ws.Range("M107").Value = 62.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 6000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2789.5557
$ws.Range("I31").Value = 1565.625
$ws.Range("K31").Value = 1565.625
$ws.Range("M31").Value = -1270.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2789.5557
$ws.Range("I34").Value = 1565.625
$ws.Range("K34").Value = 1565.625
$ws.Range("M34").Value = -1363.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 942.5
$ws.Range("I35").Value = 923.3333
$ws.Range("K35").Value = 923.3333
$ws.Range("M35").Value = -629.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2289.8333
$ws.Range("I58").Value = 1745
$ws.Range("K58").Value = 1745
$ws.Range("M58").Value = -1542

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 20874.125
$ws.Range("J60").Value = 69994
$ws.Range("L60").Value = 69994
$ws.Range("N60").Value = -71016

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7583
$ws.Range("I62").Value = 7749.5
$ws.Range("K62").Value = 7749.5
$ws.Range("M62").Value = -7125.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7583
$ws.Range("I65").Value = 7749.5
$ws.Range("K65").Value = 38747.5
$ws.Range("M65").Value = -35627.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 84999.5
$ws.Range("J87").Value = 84999.5
$ws.Range("L87").Value = 84999.5
$ws.Range("N87").Value = -87371.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 84999.5
$ws.Range("J90").Value = 84999.5
$ws.Range("L90").Value = 254998.5
$ws.Range("N90").Value = -266854.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3317.4285
$ws.Range("I99").Value = 2444.4
$ws.Range("K99").Value = 2444.4
$ws.Range("M99").Value = -946.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 717.75
$ws.Range("I105").Value = 717.75
$ws.Range("K105").Value = 717.75
$ws.Range("M105").Value = 1029.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2255.1853
$ws.Range("I122").Value = 1881.7778
$ws.Range("J122").Value = 3002
$ws.Range("K122").Value = 5645.3334
$ws.Range("L122").Value = 9006
$ws.Range("M122").Value = -3195.3334
$ws.Range("N122").Value = -13906

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3317.4285
$ws.Range("I126").Value = 2444.4
$ws.Range("K126").Value = 7333.200000000001
$ws.Range("M126").Value = -4863.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1520.9117
$ws.Range("I134").Value = 1428.0312
$ws.Range("J134").Value = 3007
$ws.Range("K134").Value = 4284.0936
$ws.Range("L134").Value = 9021
$ws.Range("M134").Value = -1749.0936
$ws.Range("N134").Value = -14091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2289.8333
$ws.Range("I136").Value = 1745
$ws.Range("K136").Value = 5235
$ws.Range("M136").Value = -2685

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 11.666667
$ws.Range("J2").Value = 9.5
$ws.Range("K2").Value = 70.00000199999999
$ws.Range("L2").Value = 57
$ws.Range("M2").Value = 42.99999800000001
$ws.Range("N2").Value = -283

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 950
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 950
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 2850
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -3196

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2310
$ws.Range("J23").Value = 1661.25
$ws.Range("L23").Value = 4983.75
$ws.Range("N23").Value = -5453.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 340
$ws.Range("I26").Value = 550
$ws.Range("K26").Value = 1650
$ws.Range("M26").Value = -1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 963.7273
$ws.Range("J38").Value = 1057
$ws.Range("L38").Value = 3171
$ws.Range("N38").Value = -3865

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3012
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3012
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 205.36363
$ws.Range("J97").Value = 229
$ws.Range("L97").Value = 687
$ws.Range("N97").Value = -1679

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 402
$ws.Range("J98").Value = 404
$ws.Range("L98").Value = 1212
$ws.Range("N98").Value = -4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 20489.938
$ws.Range("I99").Value = 15024.214
$ws.Range("K99").Value = 45072.642
$ws.Range("M99").Value = -42826.642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 36399

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 908.1613
$ws.Range("I113").Value = 1073.1
$ws.Range("J113").Value = 829.619
$ws.Range("K113").Value = 3219.3
$ws.Range("L113").Value = 2488.857
$ws.Range("M113").Value = -1049.3
$ws.Range("N113").Value = -6828.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2112.3684
$ws.Range("I132").Value = 1674.75
$ws.Range("K132").Value = 15072.75
$ws.Range("M132").Value = -12542.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2742.077
$ws.Range("I140").Value = 2513.3635
$ws.Range("K140").Value = 7540.0905
$ws.Range("M140").Value = -2360.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2880
$ws.Range("I141").Value = 2173.3333
$ws.Range("K141").Value = 6519.999899999999
$ws.Range("M141").Value = -1339.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 29999
$ws.Range("J63").Value = 29999
$ws.Range("L63").Value = 29999
$ws.Range("N63").Value = -31371

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 29999
$ws.Range("J66").Value = 29999
$ws.Range("L66").Value = 89997
$ws.Range("N66").Value = -96861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10821.25
$ws.Range("I80").Value = 12535.417
$ws.Range("J80").Value = 8250
$ws.Range("K80").Value = 12535.417
$ws.Range("L80").Value = 8250
$ws.Range("M80").Value = -11537.417
$ws.Range("N80").Value = -10246

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10821.25
$ws.Range("I83").Value = 12535.417
$ws.Range("J83").Value = 8250
$ws.Range("K83").Value = 62677.085
$ws.Range("L83").Value = 41250
$ws.Range("M83").Value = -57685.085
$ws.Range("N83").Value = -51234

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5089.5
$ws.Range("I97").Value = 1609.8462
$ws.Range("K97").Value = 1609.8462
$ws.Range("M97").Value = -1113.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3664.5881
$ws.Range("I122").Value = 3244.2424
$ws.Range("J122").Value = 4435.222
$ws.Range("K122").Value = 9732.727200000001
$ws.Range("L122").Value = 13305.666
$ws.Range("M122").Value = -7282.727200000001
$ws.Range("N122").Value = -18205.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3017.4285
$ws.Range("I126").Value = 3017.4285
$ws.Range("K126").Value = 9052.2855
$ws.Range("M126").Value = -6582.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3141.2083
$ws.Range("I132").Value = 2835.9092
$ws.Range("K132").Value = 8507.7276
$ws.Range("M132").Value = -5977.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5515.3335
$ws.Range("I7").Value = 5524.3335
$ws.Range("J7").Value = 5479.3335
$ws.Range("K7").Value = 5524.3335
$ws.Range("L7").Value = 5479.3335
$ws.Range("M7").Value = -5412.3335
$ws.Range("N7").Value = -5703.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3021.9143
$ws.Range("I22").Value = 2825.4
$ws.Range("J22").Value = 3169.3
$ws.Range("K22").Value = 2825.4
$ws.Range("L22").Value = 3169.3
$ws.Range("M22").Value = -2530.4
$ws.Range("N22").Value = -3759.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3021.9143
$ws.Range("I27").Value = 2825.4
$ws.Range("J27").Value = 3169.3
$ws.Range("K27").Value = 2825.4
$ws.Range("L27").Value = 3169.3
$ws.Range("M27").Value = -2718.4
$ws.Range("N27").Value = -3383.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2440.6667
$ws.Range("I46").Value = 1297
$ws.Range("J46").Value = 3203.111
$ws.Range("K46").Value = 1297
$ws.Range("L46").Value = 3203.111
$ws.Range("M46").Value = -1109
$ws.Range("N46").Value = -3579.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5837.0835
$ws.Range("I93").Value = 4859
$ws.Range("J93").Value = 6535.7144
$ws.Range("K93").Value = 4859
$ws.Range("L93").Value = 6535.7144
$ws.Range("M93").Value = -3611
$ws.Range("N93").Value = -9031.714400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5843.7334
$ws.Range("I122").Value = 5169
$ws.Range("J122").Value = 7699.25
$ws.Range("K122").Value = 15507
$ws.Range("L122").Value = 23097.75
$ws.Range("M122").Value = -13057
$ws.Range("N122").Value = -27997.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5515.3335
$ws.Range("I126").Value = 5524.3335
$ws.Range("J126").Value = 5479.3335
$ws.Range("K126").Value = 16573.0005
$ws.Range("L126").Value = 16438.0005
$ws.Range("M126").Value = -14103.0005
$ws.Range("N126").Value = -21378.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3754.543
$ws.Range("I132").Value = 3837.2778
$ws.Range("J132").Value = 3666.9412
$ws.Range("K132").Value = 11511.8334
$ws.Range("L132").Value = 11000.8236
$ws.Range("M132").Value = -8981.8334
$ws.Range("N132").Value = -16060.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3681.26
$ws.Range("I136").Value = 4203.75
$ws.Range("K136").Value = 12611.25
$ws.Range("M136").Value = -10061.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1568.7142
$ws.Range("I81").Value = 970
$ws.Range("K81").Value = 1940
$ws.Range("M81").Value = -879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1568.7142
$ws.Range("I84").Value = 970
$ws.Range("K84").Value = 9700
$ws.Range("M84").Value = -4396

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4145.5654
$ws.Range("I96").Value = 4679.2666
$ws.Range("J96").Value = 3144.875
$ws.Range("K96").Value = 4679.2666
$ws.Range("L96").Value = 3144.875
$ws.Range("M96").Value = -3306.2666
$ws.Range("N96").Value = -5890.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 163333
$ws.Range("J116").Value = 163333
$ws.Range("L116").Value = 163333
$ws.Range("N116").Value = -172511

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 84924.5
$ws.Range("J119").Value = 84924.5
$ws.Range("L119").Value = 84924.5
$ws.Range("N119").Value = -94600.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2756.9268
$ws.Range("I132").Value = 2104.2942
$ws.Range("J132").Value = 5926.857
$ws.Range("K132").Value = 6312.882599999999
$ws.Range("L132").Value = 17780.571
$ws.Range("M132").Value = -3782.882599999999
$ws.Range("N132").Value = -22840.571
